$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PIR")

$row = 17

# Column A holds a date-formatted string ("2026-02-01"); Excel would
# normally auto-convert that to a date serial number. Force it to stay
# plain text, then clear the formatting so no extra cell style sticks.
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2026-02-01"
$cellA.ClearFormats()

$ws.Cells.Item($row, 2).Value = "13:43:58"
$ws.Cells.Item($row, 3).Value = "13:00"
$ws.Cells.Item($row, 4).Value = "Bathroom"
$ws.Cells.Item($row, 5).Value = "No Motion"
$ws.Cells.Item($row, 6).Value = "Inactive"
